# Adds a new "Player Info" sheet before "ODI Batting", and updates the
# MATCH_CARD_LINK columns in "ODI Batting" / "ODI Bowling" to a new
# MATCH_CODE column containing just the numeric match code instead of
# the full scorecard URL.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet before "ODI Batting"
# ---------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

# Re-fetch the other sheets by name now that a new sheet has been
# inserted - worksheet references obtained before an Add() can become
# stale once the sheet collection shifts.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $infoHeaders.Length; $col++) {
    $infoSheet.Cells.Item(1, $col).Value = $infoHeaders[$col - 1]
}

$infoSheet.Cells.Item(2, 1).NumberFormat = "@"
$infoSheet.Cells.Item(2, 1).Value = "3978"
$infoSheet.Cells.Item(2, 2).Value = "Todd Duncan Astle"
$infoSheet.Cells.Item(2, 3).Value = "Right Handed"
$infoSheet.Cells.Item(2, 4).Value = "Right Arm Leg Break"

# Match the bold / bordered / centered header style used elsewhere in
# the workbook.
$infoHeaderRange = $infoSheet.Range("A1:D1")
$infoHeaderRange.Font.Bold = $true
$infoHeaderRange.Borders.LineStyle = 1
$infoHeaderRange.HorizontalAlignment = -4108
$infoHeaderRange.VerticalAlignment = -4160

$infoSheet.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# Helper: rewrite a MATCH_CARD_LINK column (full scorecard URL) into a
# MATCH_CODE column (just the numeric match code).
# ---------------------------------------------------------------------
function Update-MatchCodeColumn($sheet, $colLetter, $lastRow) {
    $headerCell = $sheet.Range($colLetter + "1")
    $headerCell.Value = "MATCH_CODE"

    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $sheet.Range($colLetter + $row)
        $link = $cell.Value2
        if ($link) {
            $code = $link.ToString().Split("=")[-1]
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }
}

# ---------------------------------------------------------------------
# 2. "ODI Batting": column D (MATCH_CARD_LINK -> MATCH_CODE)
# ---------------------------------------------------------------------
Update-MatchCodeColumn $battingSheet "D" 10

# ---------------------------------------------------------------------
# 3. "ODI Bowling": column B (MATCH_CARD_LINK -> MATCH_CODE)
# ---------------------------------------------------------------------
Update-MatchCodeColumn $bowlingSheet "B" 9
